# Commit: "Change names from *img to img*"
# Rename the seven "<letter>img" sheets to "img<letter>" and leave the
# last one ("imge", formerly "eimg") selected as the active tab.

$wb = $excel.ActiveWorkbook

$renames = @(
    @("himg", "imgh"),
    @("timg", "imgt"),
    @("simg", "imgs"),
    @("gimg", "imgg"),
    @("wimg", "imgw"),
    @("bimg", "imgb"),
    @("eimg", "imge")
)

foreach ($pair in $renames) {
    $oldName = $pair[0]
    $newName = $pair[1]
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $newName
}

# Select/activate the renamed "imge" sheet (was "eimg") as the active tab.
$wb.Worksheets.Item("imge").Activate()
